# Weekly Summary refresh: new reporting period (Jan 01 -> Jan 11, 2026),
# refreshed metrics, two new vendor rows, renumbered sections, and a new
# "Top canceller" recommendation line.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Weekly Summary")

# ---------------------------------------------------------------------
# 1. Header title / reporting period
# ---------------------------------------------------------------------
$ws.Range("A1").Value = '📊 Weekly Summary • Jan 01 → Jan 11, 2026'

# ---------------------------------------------------------------------
# 2. USERS block
# ---------------------------------------------------------------------
$ws.Range("B4").Value = 278
$ws.Range("B5").Value = 74
$ws.Range("B6").Value = 278

# ---------------------------------------------------------------------
# 3. ORDERS block
# ---------------------------------------------------------------------
$ws.Range("B9").Value = 186
$ws.Range("D9").Value = 147
$ws.Range("F9").Value = 38
$ws.Range("B10").Value = 79

# ---------------------------------------------------------------------
# 4. FINANCIALS block
# ---------------------------------------------------------------------
$ws.Range("B13").Value = 28455
$ws.Range("D13").Value = 1980
$ws.Range("F13").Value = 30435

# ---------------------------------------------------------------------
# 5. TOP VENDORS table - update existing 4 rows (17-20)
# ---------------------------------------------------------------------
$ws.Range("A17").Value = 'ክርስቲና ምግብ ቤት'
$ws.Range("D17").Value = 2

$ws.Range("A18").Value = 'Test Vendor'

$ws.Range("A19").Value = 'መቅዲ ምግብ ቤት'

$ws.Range("A20").Value = 'ቲጂ አቡዳቢ #5kilo'
$ws.Range("B20").Value = 12930
$ws.Range("C20").Value = 14060
$ws.Range("D20").Value = 15

# ---------------------------------------------------------------------
# 6. TOP VENDORS table - insert two new vendor rows (21, 22), pushing
#    everything below down by two rows. Insert() copies the formatting
#    (number formats / styles) from the row above, matching rows 17-20.
# ---------------------------------------------------------------------
$ws.Rows.Item(21).Insert()
$ws.Rows.Item(21).Insert()

$ws.Range("A21").Value = 'ቲጂ አቡዳቢ #6kilo'
$ws.Range("B21").Value = 8680
$ws.Range("C21").Value = 9285
$ws.Range("D21").Value = 13

$ws.Range("A22").Value = 'ጤና ምግብ ቤት'
$ws.Range("B22").Value = 6845
$ws.Range("C22").Value = 7330
$ws.Range("D22").Value = 7

# ---------------------------------------------------------------------
# 7. TOP MEAL / TOP CAMPUS (now at rows 25/26 after the insert above)
# ---------------------------------------------------------------------
$ws.Range("B25").Value = 'ሙሉ ኮርኒስ ×37'
$ws.Range("B26").Value = '6kilo (84 orders)'

# ---------------------------------------------------------------------
# 8. DELIVERY SQUAD summary row (now at row 30)
# ---------------------------------------------------------------------
$ws.Range("A30").Value = 2
$ws.Range("B30").Value = 15

# ---------------------------------------------------------------------
# 9. RECOMMENDATIONS section (now rows 36-38), plus a brand new row 39
#    for the "Top canceller" insight. Insert a row after row 38 (the
#    current last recommendation line) so the trailing blank row /
#    footer shift down to match the diff.
# ---------------------------------------------------------------------
$ws.Range("A36").Value = 'High cancellation rate this week (20.4%). Investigate top cancelled meals/vendors and contact them.'
$ws.Range("A37").Value = 'Vendors with low reliability: ክርስቲና ምግብ ቤት, Test Vendor, መቅዲ ምግብ ቤት. Consider warnings, training, or temporary delisting.'
# A38 ("DG acceptance low ...") is unchanged.

$ws.Rows.Item(39).Insert()
$ws.Range("A39").Value = 'Top canceller: ቲጂ አቡዳቢ #5kilo with 15 cancels. Investigate immediately.'

# ---------------------------------------------------------------------
# 10. Footer timestamp (now at E41)
# ---------------------------------------------------------------------
$ws.Range("E41").Value = 'Prepared on: 2026-01-11 11:06 UTC'

# ---------------------------------------------------------------------
# 11. Column widths: B 18 -> 19, F 7 -> 9 (Excel stores width as the
#     ColumnWidth plus a constant ~0.8333 padding offset, so subtract
#     that offset from the target display width before assigning).
# ---------------------------------------------------------------------
$padding = 5.0 / 6.0
$ws.Columns.Item(2).ColumnWidth = 19 - $padding
$ws.Columns.Item(6).ColumnWidth = 9 - $padding

# ---------------------------------------------------------------------
# 12. Chart: extend the category/value series ranges from rows 17:20 to
#     17:22 to include the two new vendors, and move the chart anchor
#     down by two rows (2 * default 15pt row height = 30pt) so it keeps
#     sitting just below the (now longer) vendor table.
# ---------------------------------------------------------------------
$co = $ws.ChartObjects().Item(1)
$co.Chart.SeriesCollection(1).Formula = "=SERIES('Weekly Summary'!`$B`$16,'Weekly Summary'!`$A`$17:`$A`$22,'Weekly Summary'!`$B`$17:`$B`$22,1)"
$co.Top = $co.Top + 30
